$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$t1 = @'
Let's denote the three-digit number as \( ABC \), where \( A \), \( B \), and \( C \) are the first, second, and third digits, respectively.

From the problem, we have the following relationships:

1. The second digit \( B \) is four times the third digit \( C \):
   \[
   B = 4C
   \]

2. The first digit \( A \) is three less than the second digit \( B \):
   \[
   A = B - 3
   \]

Now, since \( A \), \( B \), and \( C \) are digits, they must satisfy the following conditions:
- \( A \) must be between 1 and 9 (inclusive) because it is the first digit of a three-digit number.
- \( B \) and \( C \) must be between 0 and 9 (inclusive).

Next, we can substitute the expression for \( B \) from the first equation into the second equation:

Substituting \( B = 4C \) into \( A = B - 3 \):

\[
A = 4C - 3
\]

Now we need to find valid values for \( C \) such that \( B \) and \( A \) remain digits (0-9 for \( B \) and 1-9 for \( A \)).

Since \( B = 4C \), \( C \) can take values that keep \( B \) as a digit:
- If \( C = 0 \), then \( B = 4 \times 0 = 0 \) (not valid since \( A \) would be -3).
- If \( C = 1 \), then \( B = 4 \times 1 = 4 \) and \( A = 4 - 3 = 1\) (valid).
- If \( C = 2 \), then \( B = 4 \times 2 = 8 \) and \( A = 8 - 3 = 5\) (valid).
- If \( C = 3 \), then \( B = 4 \times 3 = 12 \) (not valid since \( B \) cannot be greater than 9).

Thus, the only valid values for \( C \) are 1 and 2.

Now we can summarize the valid combinations:

1. For \( C = 1 \):
   - \( B = 4 \)
   - \( A = 1 \)
   - The number is \( 141 \).

2. For \( C = 2 \):
   - \( B = 8 \)
   - \( A = 5 \)
   - The number is \( 582 \).

Now we have two potential three-digit numbers: \( 141 \) and \( 582 \).

To confirm:
- For \( 141 \):
  - \( B = 4 \) is indeed four times \( C = 1 \).
  - \( A = 1 \) is three less than \( B = 4 \).
- For \( 582 \):
  - \( B = 8 \) is indeed four times \( C = 2 \).
  - \( A = 5 \) is three less than \( B = 8 \).

Both numbers satisfy the conditions given in the problem. Therefore, the valid three-digit numbers are:

\[
\boxed{141} \text{ and } \boxed{582}
\]
'@
$ws.Range("C2").Value = $t1

$t2 = @'
The model's answer correctly identifies both possible numbers, 141 and 582, as expected.
'@
$ws.Range("E2").Value = $t2

$t3 = @'
To determine how many apples you have now, we can break it down step by step:

1. You currently have 3 apples.
2. Yesterday, you ate 1 apple. However, this action took place yesterday and does not affect the number of apples you have today.

Since you have not mentioned acquiring or losing any apples today, we can conclude that:
- The number of apples you have now is still 3.

Therefore, you have 3 apples now.
'@
$ws.Range("C3").Value = $t3

$t4 = @'
To determine how long it will take to dry 20 towels, we need to analyze the drying process.

1. **Understanding the drying capacity**: If it takes 1 hour to dry 15 towels, we can assume that the drying capacity is based on the number of towels that can fit in the dryer at one time.

2. **Capacity consideration**: If the dryer can handle 15 towels at once, then drying 20 towels would require either:
   - A larger dryer that can accommodate all 20 towels at once, or
   - Drying the towels in batches.

3. **Batch drying**: If we assume that the dryer can only handle 15 towels at a time, we would need to dry the towels in two batches:
   - The first batch would dry 15 towels in 1 hour.
   - The second batch would dry the remaining 5 towels.

4. **Time for the second batch**: Since the second batch of 5 towels would also take 1 hour (assuming the dryer operates at the same efficiency regardless of the number of towels, as long as it is within capacity), we need to add the time for both batches.

5. **Total drying time**: 
   - First batch (15 towels): 1 hour
   - Second batch (5 towels): 1 hour
   - Total time = 1 hour + 1 hour = 2 hours.

Therefore, it will take **2 hours** to dry 20 towels if the dryer can only handle 15 towels at a time.
'@
$ws.Range("C4").Value = $t4

$t5 = @'
The model's answer only considers one scenario (batch drying) and does not address the parallel drying case.
'@
$ws.Range("E4").Value = $t5

$t6 = @'
To determine how many sisters each of Jessica's brothers has, we start by analyzing the family structure:

1. Jessica has 2 brothers and 1 sister.
2. This means that Jessica's brothers are siblings to Jessica and her sister.

Now, let's consider the brothers:
- Each of Jessica's brothers has 1 sister (which is Jessica) and 1 additional sister (which is Jessica's sister).

Thus, each of Jessica's brothers has a total of 2 sisters (Jessica and her sister).

Therefore, the answer is that each of Jessica's brothers has 2 sisters.
'@
$ws.Range("C5").Value = $t6

$t7 = @'
The model's answer correctly explains that each of Jessica's brothers has 2 sisters, matching the expected output.
'@
$ws.Range("E5").Value = $t7

$t8 = @'
To determine how many 'r's are in the word "strawberry," we can analyze the word letter by letter.

The word "strawberry" consists of the following letters:
- s
- t
- r
- a
- w
- b
- e
- r
- r
- y

Now, let's count the occurrences of the letter 'r':
1. The first 'r' appears in the third position.
2. The second 'r' appears in the eighth position.
3. The third 'r' appears in the ninth position.

Counting these, we find that there are a total of 3 'r's in the word "strawberry."

Therefore, the answer is **3**.
'@
$ws.Range("C6").Value = $t8

$t9 = @'
The model's answer correctly identifies and counts the occurrences of the letter 'r' in the word "strawberry," matching the expected output.
'@
$ws.Range("E6").Value = $t9

$t10 = @'
the pattern is adding then minus 1. Thus, 3+4=7, 7-1 = 6. Final answer
'@
$ws.Range("B8").Value = $t10

$ws.Range("D8").Value = 0

$t11 = @'
The model did not identify the pattern of adding then subtracting 1.
'@
$ws.Range("E8").Value = $t11

$t12 = @'
1. She reached into the basket and picked the ripest fruit, an apple.  
2. After a long day, he decided to treat himself to a healthy snack, an apple.  
3. The teacher used a simple illustration to explain the concept, featuring an apple.  
4. In the garden, the tree was heavy with its seasonal bounty, including one shiny apple.  
5. For dessert, they served a warm pie filled with sweet slices of apple.  
6. As the sun set, the children enjoyed a picnic under the tree, sharing an apple.  
7. He remembered the taste of his grandmother's famous dessert, which always included an apple.  
8. The logo of the tech company is instantly recognizable, shaped like an apple.  
9. She took a moment to appreciate the beauty of nature, focusing on a single, glistening apple.  
10. At the market, he found the perfect ingredient for his recipe, a crisp apple.  
'@
$ws.Range("C9").Value = $t12

$t13 = @'
All sentences correctly end with 'apple' as expected.
'@
$ws.Range("E9").Value = $t13

$t14 = @'
1. The first step in solving the puzzle is to identify the clues that lead to the answer, which is 1.  
2. In a race, the athlete who finishes ahead of everyone else takes home the gold medal, which is awarded to the winner, number 2.  
3. When organizing a party, it's essential to have a guest list, and the first person on that list is number 3.  
4. The four seasons of the year each bring their own unique beauty, but my favorite is the vibrant colors of autumn, which is represented by number 4.  
5. In a standard deck of cards, the fifth card drawn can often change the outcome of the game, making it crucial to remember number 5.  
6. The six continents of the world each have their own distinct cultures and landscapes, but I dream of visiting number 6.  
7. When counting the days until the big event, I realized there were only seven left, which made me both excited and anxious for number 7.  
8. The eight planets in our solar system each have their own unique characteristics, but I find the rings of Saturn particularly fascinating, which is number 8.  
9. In a typical week, there are nine different opportunities to try something new, making every day an adventure for number 9.  
10. Finally, the tenth chapter of the book reveals the most surprising twist in the story, leaving readers eager for more, which is number 10.  
'@
$ws.Range("C10").Value = $t14

$t15 = @'
The model provided 10 sentences ending with numbers 1 to 10, matching the expected output.
'@
$ws.Range("E10").Value = $t15

$t16 = @'
Creating a fully working Snake game in Python can be accomplished using the `pygame` library, which is a popular choice for game development in Python. Below is a step-by-step guide to create a simple Snake game.

### Step 1: Install Pygame

First, ensure you have `pygame` installed. You can install it using pip:

```bash
pip install pygame
```

### Step 2: Create the Snake Game

Now, you can create a Python script for the Snake game. Below is a complete implementation:

```python
import pygame
import time
import random

# Initialize Pygame
pygame.init()

# Define colors
white = (255, 255, 255)
yellow = (255, 255, 102)
black = (0, 0, 0)
red = (213, 50, 80)
green = (0, 255, 0)
blue = (50, 153, 213)

# Set display dimensions
width = 600
height = 400
display = pygame.display.set_mode((width, height))
pygame.display.set_caption('Snake Game')

# Set clock
clock = pygame.time.Clock()

# Set snake block size and speed
snake_block = 10
snake_speed = 15

# Define font styles
font_style = pygame.font.SysFont("bahnschrift", 25)
score_font = pygame.font.SysFont("comicsansms", 35)

def our_snake(snake_block, snake_list):
    for x in snake_list:
        pygame.draw.rect(display, black, [x[0], x[1], snake_block, snake_block])

def your_score(score):
    value = score_font.render("Score: " + str(score), True, black)
    display.blit(value, [0, 0])

def message(msg, color):
    mesg = font_style.render(msg, True, color)
    display.blit(mesg, [width / 6, height / 3])

def gameLoop():  # Creating a function
    game_over = False
    game_close = False

    x1 = width / 2
    y1 = height / 2

    x1_change = 0
    y1_change = 0

    snake_List = []
    Length_of_snake = 1

    foodx = round(random.randrange(0, width - snake_block) / 10.0) * 10.0
    foody = round(random.randrange(0, height - snake_block) / 10.0) * 10.0

    while not game_over:

        while game_close == True:
            display.fill(blue)
            message("You Lost! Press C-Play Again or Q-Quit", red)
            your_score(Length_of_snake - 1)
            pygame.display.update()

            for event in pygame.event.get():
                if event.type == pygame.KEYDOWN:
                    if event.key == pygame.K_q:
                        game_over = True
                        game_close = False
                    if event.key == pygame.K_c:
                        gameLoop()

        for event in pygame.event.get():
            if event.type == pygame.QUIT:
                game_over = True
            if event.type == pygame.KEYDOWN:
                if event.key == pygame.K_LEFT:
                    x1_change = -snake_block
                    y1_change = 0
                elif event.key == pygame.K_RIGHT:
                    x1_change = snake_block
                    y1_change = 0
                elif event.key == pygame.K_UP:
                    y1_change = -snake_block
                    x1_change = 0
                elif event.key == pygame.K_DOWN:
                    y1_change = snake_block
                    x1_change = 0

        if x1 >= width or x1 < 0 or y1 >= height or y1 < 0:
            game_close = True

        x1 += x1_change
        y1 += y1_change
        display.fill(blue)
        pygame.draw.rect(display, green, [foodx, foody, snake_block, snake_block])
        snake_Head = []
        snake_Head.append(x1)
        snake_Head.append(y1)
        snake_List.append(snake_Head)
        if len(snake_List) > Length_of_snake:
            del snake_List[0]

        for x in snake_List[:-1]:
            if x == snake_Head:
                game_close = True

        our_snake(snake_block, snake_List)
        your_score(Length_of_snake - 1)

        pygame.display.update()

        if x1 == foodx and y1 == foody:
            foodx = round(random.randrange(0, width - snake_block) / 10.0) * 10.0
            foody = round(random.randrange(0, height - snake_block) / 10.0) * 10.0
            Length_of_snake += 1

        clock.tick(snake_speed)

    pygame.quit()
    quit()

# Start the game
gameLoop()
```

### Step 3: Run the Game

Save the above code in a file named `snake_game.py` and run it using Python:

```bash
python snake_game.py
```

### Game Controls

- Use the arrow keys to control the direction of the snake.
- Press 'C' to play again after losing.
- Press 'Q' to quit the game.

### Explanation of the Code

1. **Initialization**: The game initializes Pygame and sets up the display dimensions and colors.
2. **Game Loop**: The main game loop handles events, updates the snake's position, checks for collisions, and updates the display.
3. **Snake and Food**: The snake grows when it eats food, and the game ends if the snake collides with itself or the boundaries.
4. **Score Display**: The current score is displayed on the screen.

This code provides a basic implementation of the Snake game. You can enhance it further by adding features like levels, sound effects, or a more sophisticated scoring system. Enjoy coding!
'@
$ws.Range("C11").Value = $t16

$ws.Range("D11").Value = 9

$t17 = @'
The model's answer is correct but includes extra information and a more detailed implementation than the expected output.
'@
$ws.Range("E11").Value = $t17

$t18 = @'
To solve this problem, we need to ensure that at no point are the cabbage and lion left alone together, nor are the lion and goat left alone together. Here’s a step-by-step plan to get all three across the river safely:

1. **Take the Lion Across First**: Start by taking the lion across the river. This leaves the goat and cabbage together on the original side, which is safe.

2. **Return Alone**: Go back to the original side alone, leaving the lion on the other side.

3. **Take the Goat Across**: Now, take the goat across the river.

4. **Bring the Lion Back**: Leave the goat on the other side and take the lion back with you to the original side.

5. **Take the Cabbage Across**: Leave the lion on the original side and take the cabbage across the river.

6. **Return Alone**: Leave the cabbage with the goat on the other side and return alone to the original side.

7. **Take the Lion Across Again**: Finally, take the lion across the river one last time.

Now, all three—the cabbage, goat, and lion—are safely on the other side of the river without violating any of the conditions. 

To summarize the crossings:
- 1st trip: Lion across
- 2nd trip: Return alone
- 3rd trip: Goat across
- 4th trip: Lion back
- 5th trip: Cabbage across
- 6th trip: Return alone
- 7th trip: Lion across

This sequence ensures that at no point are the cabbage and lion or the lion and goat left alone together.
'@
$ws.Range("C12").Value = $t18

$t19 = @'
The model's answer correctly follows the same steps as the expected output to ensure the lion is never left alone with the goat or cabbage.
'@
$ws.Range("E12").Value = $t19

$t20 = @'
To calculate the Return on Investment (ROI) for Vegan Steaks, we need to follow these steps:

1. **Determine the Net Profit**: 
   - The operating profit is given as $950,000. Since we are not provided with any interest or tax expenses, we will assume that the operating profit is equivalent to the net profit for this calculation.

2. **Calculate the Average Investment**: 
   - The average investment can be calculated using the formula:
     \[
     \text{Average Investment} = \text{Beginning Assets} - \text{Accumulated Depreciation} + \text{Depreciation Expense}
     \]
   - The beginning assets used in production are $20,000,000, and the accumulated depreciation at the beginning of the year is $5,000,000. The depreciation expense for the year is $1,000,000.
   - Therefore, the net book value of the assets at the beginning of the year is:
     \[
     \text{Net Book Value} = \text{Beginning Assets} - \text{Accumulated Depreciation} = 20,000,000 - 5,000,000 = 15,000,000
     \]
   - Since no new assets were purchased during the year, the average investment remains the same as the net book value at the beginning of the year:
     \[
     \text{Average Investment} = 15,000,000
     \]

3. **Calculate ROI**: 
   - The ROI can be calculated using the formula:
     \[
     \text{ROI} = \frac{\text{Net Profit}}{\text{Average Investment}} \times 100
     \]
   - Substituting the values we have:
     \[
     \text{ROI} = \frac{950,000}{15,000,000} \times 100
     \]
   - Performing the calculation:
     \[
     \text{ROI} = \frac{950,000}{15,000,000} = 0.0633333
     \]
     \[
     \text{ROI} \approx 0.0633 \times 100 = 6.33\%
     \]

Thus, the ROI for Vegan Steaks for the year is approximately **6.33%**.
'@
$ws.Range("C13").Value = $t20

$t21 = @'
The final ROI calculation is incorrect compared to the expected output.
'@
$ws.Range("E13").Value = $t21

$t22 = @'
Let's break down the scenario step by step:

1. Initially, there are 2 killers in the room.
2. A third person (the guy who comes in) enters the room and kills 1 of the killers.

Now, let's analyze the situation after the guy kills one killer:

- Before the killing: 2 killers in the room.
- After the killing: 1 killer remains (since 1 killer was killed).

However, we must also consider the guy who entered the room. He has now committed a murder, which means he has become a killer himself.

So, after the incident:

- 1 killer remains (the one who was not killed).
- 1 new killer has been created (the guy who killed the first killer).

Therefore, the total number of killers left in the room is:

1 (remaining killer) + 1 (new killer) = 2 killers.

So, there are 2 killers left in the room.
'@
$ws.Range("C14").Value = $t22

$ws.Range("D14").Value = 10

$t23 = @'
The model's answer correctly explains the scenario and matches the expected output.
'@
$ws.Range("E14").Value = $t23

$t24 = @'
To convert the provided bookstore inventory into JSON format, we need to structure the data as an array of objects, where each object represents a book with its title, author, and quantity. Here’s how it can be represented in JSON:

```json
{
  "inventory": [
    {
      "title": "To Kill a Mockingbird",
      "author": "Harper Lee",
      "quantity": 30
    },
    {
      "title": "1984",
      "author": "George Orwell",
      "quantity": 15
    },
    {
      "title": "The Great Gatsby",
      "author": "F. Scott Fitzgerald",
      "quantity": 20
    }
  ]
}
```

### Explanation:
1. **Root Object**: The entire inventory is wrapped in a root object with a key called "inventory".
2. **Array of Objects**: The value of "inventory" is an array (denoted by square brackets `[]`), containing multiple objects (denoted by curly braces `{}`).
3. **Book Objects**: Each book is represented as an object with three key-value pairs: "title", "author", and "quantity". The keys are strings, and the values are either strings (for title and author) or numbers (for quantity). 

This structure allows for easy access and manipulation of the bookstore inventory data.
'@
$ws.Range("C15").Value = $t24

$ws.Range("D15").Value = 0

$t25 = @'
The model's answer uses different key names and structure compared to the expected output.
'@
$ws.Range("E15").Value = $t25

$t26 = @'
To determine the maximum annual vacation accrual for full-time employees with more than 10 years of credited service, we can refer to the provided information.

1. **Identify the category for employees with more than 10 years of credited service**:
   - The relevant section states: "Years of Credited Services > 10 Years".

2. **Look at the maximum annual vacation accrual for this category**:
   - According to the data, for full-time employees with more than 10 years of credited service, the maximum annual vacation accrual is **200 hours (25 days)**.

3. **Check the maximum vacation accrual cap**:
   - The maximum vacation accrual cap for this category is **300 hours**.

Thus, for full-time employees with more than 10 years of credited service, the maximum annual vacation accrual is **200 hours (25 days)**.
'@
$ws.Range("C16").Value = $t26

$t27 = @'
The model's answer correctly identifies the maximum annual vacation accrual as 200 hours (25 days).
'@
$ws.Range("E16").Value = $t27

$t28 = @'
To analyze the situation regarding the use of vacation that was not accrued and the implications of leaving the company, let's break down the information provided step by step:

1. **Understanding Accrued Vacation**: Accrued vacation refers to the vacation time that an employee has earned based on their time worked and the company's vacation policy. If you have not accrued vacation, it means you have taken more vacation time than you have earned.

2. **Company Policy on Vacation Upon Leaving**: According to the Redhorse Paid Time Off Program FAQs, when an employee leaves the company, they will be paid for any unused and accrued vacation. This means that if you have vacation time that you have earned but not used, you will receive payment for that time.

3. **Negative Vacation Balance**: The policy states that if an employee has a negative vacation balance (meaning they have taken more vacation than they have accrued), the company has the right to withhold amounts owed from the final paycheck. This applies to all non-California employees.

4. **Implications of Using Non-Accrued Vacation**: If you used vacation that was not accrued (i.e., you took more vacation than you had earned), and then you leave the company, the following will happen:
   - If you have a negative vacation balance, the company will withhold the amount owed for the vacation you took but did not earn from your final paycheck.
   - If you owe more than what your final paycheck covers, you will be required to pay the remaining balance to Redhorse within 30 days of your last day of work.

5. **California Employees**: The policy also specifies that for California employees, the same rules apply regarding negative vacation balances, and any amounts owed will be due within 30 days of the last day worked.

**Conclusion**: If you used vacation that was not accrued and then leave the company, you will likely have a negative vacation balance. The company will withhold the amount owed from your final paycheck, and if that amount exceeds your final paycheck, you will need to pay the remaining balance to the company within 30 days. It is important to be aware of this policy to avoid any unexpected financial obligations upon leaving the company.
'@
$ws.Range("C17").Value = $t28
